$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.612942462218352
$ws.Range("D2").Value = 8.083804492557505
$ws.Range("E2").Value = 12.77847057985256
$ws.Range("F2").Value = 36.4946390875902
$ws.Range("G2").Value = 41.0341322947306
$ws.Range("H2").Value = 17.3313590401117
$ws.Range("I2").Value = 25.60809687858471
$ws.Range("J2").Value = 9.944838780935042
$ws.Range("K2").Value = 15.9145656672906
$ws.Range("M2").Value = 17.955369247789
$ws.Range("C3").Value = 4.552621506933073
$ws.Range("D3").Value = 8.064612036372964
$ws.Range("E3").Value = 12.79880000593528
$ws.Range("F3").Value = 36.60589574905942
$ws.Range("G3").Value = 41.1566065400751
$ws.Range("H3").Value = 17.40635407738986
$ws.Range("I3").Value = 25.69531232492648
$ws.Range("J3").Value = 9.975467833173036
$ws.Range("K3").Value = 15.44305825294019
$ws.Range("M3").Value = 17.770164642657
$ws.Range("C4").Value = 4.514821498750857
$ws.Range("D4").Value = 8.053485507296664
$ws.Range("E4").Value = 12.81317503421434
$ws.Range("F4").Value = 36.68547063445276
$ws.Range("G4").Value = 41.2485431951707
$ws.Range("H4").Value = 17.45641433328324
$ws.Range("I4").Value = 25.75663094191937
$ws.Range("J4").Value = 9.995470777719611
$ws.Range("K4").Value = 15.1480008243316
$ws.Range("M4").Value = 17.65809833795838
$ws.Range("C5").Value = 4.499235483785327
$ws.Range("D5").Value = 8.049119075677567
$ws.Range("E5").Value = 12.81950889697689
$ws.Range("F5").Value = 36.720715838463
$ws.Range("G5").Value = 41.29018601417958
$ws.Range("H5").Value = 17.47782066553027
$ws.Range("I5").Value = 25.78356221524528
$ws.Range("J5").Value = 10.00392346337046
$ws.Range("K5").Value = 15.02654370509079
$ws.Range("M5").Value = 17.6128882148308
$ws.Range("C6").Value = 4.496636717536355
$ws.Range("D6").Value = 8.048404229086829
$ws.Range("E6").Value = 12.82058937330102
$ws.Range("F6").Value = 36.72673803108458
$ws.Range("G6").Value = 41.29735216674604
$ws.Range("H6").Value = 17.48143587320762
$ws.Range("I6").Value = 25.78815121451462
$ws.Range("J6").Value = 10.00534523740669
$ws.Range("K6").Value = 15.00630763613475
$ws.Range("M6").Value = 17.60540999583824
$ws.Range("C7").Value = 4.514612025229956
$ws.Range("D7").Value = 8.053425937964558
$ws.Range("E7").Value = 12.81325852796201
$ws.Range("F7").Value = 36.68593457409587
$ws.Range("G7").Value = 41.24908793102931
$ws.Range("H7").Value = 17.45669895546511
$ws.Range("I7").Value = 25.75698629102344
$ws.Range("J7").Value = 9.995583552975496
$ws.Range("K7").Value = 15.14636750349279
$ws.Range("M7").Value = 17.65748670875612
$ws.Range("C8").Value = 4.592307767990842
$ws.Range("D8").Value = 8.077051848738455
$ws.Range("E8").Value = 12.78508735064172
$ws.Range("F8").Value = 36.53065438885918
$ws.Range("G8").Value = 41.07286979629129
$ws.Range("H8").Value = 17.35638279228329
$ws.Range("I8").Value = 25.63655070334837
$ws.Range("J8").Value = 9.955151587603755
$ws.Range("K8").Value = 15.75324399747029
$ws.Range("M8").Value = 17.8911941795647
$ws.Range("C9").Value = 4.738210070220929
$ws.Range("D9").Value = 8.128495739926306
$ws.Range("E9").Value = 12.7448605613856
$ws.Range("F9").Value = 36.31612245361448
$ws.Range("G9").Value = 40.86137871086675
$ws.Range("H9").Value = 17.19162628952993
$ws.Range("I9").Value = 25.4624354709476
$ws.Range("J9").Value = 9.885339371930337
$ws.Range("K9").Value = 16.89216666656584
$ws.Range("M9").Value = 18.3606725432843
$ws.Range("C10").Value = 4.840943009719735
$ws.Range("D10").Value = 8.169264149286295
$ws.Range("E10").Value = 12.72445788788957
$ws.Range("F10").Value = 36.21413088044881
$ws.Range("G10").Value = 40.78930838835527
$ws.Range("H10").Value = 17.09023172885264
$ws.Range("I10").Value = 25.37289533488156
$ws.Range("J10").Value = 9.839797198535418
$ws.Range("K10").Value = 17.68912136660019
$ws.Range("M10").Value = 18.709825257865
$ws.Range("C11").Value = 4.886605197147892
$ws.Range("D11").Value = 8.188423605501026
$ws.Range("E11").Value = 12.71716172007372
$ws.Range("F11").Value = 36.179957323652
$ws.Range("G11").Value = 40.77489541397744
$ws.Range("H11").Value = 17.0484076653595
$ws.Range("I11").Value = 25.34060262449167
$ws.Range("J11").Value = 9.820321552566803
$ws.Range("K11").Value = 18.04146817035956
$ws.Range("M11").Value = 18.86903072624208
$ws.Range("C12").Value = 4.903733807403236
$ws.Range("D12").Value = 8.195764106156254
$ws.Range("E12").Value = 12.71468403160888
$ws.Range("F12").Value = 36.16878468249028
$ws.Range("G12").Value = 40.77209827235728
$ws.Range("H12").Value = 17.03319141893156
$ws.Range("I12").Value = 25.32959568209164
$ws.Range("J12").Value = 9.813124788506954
$ws.Range("K12").Value = 18.17331366089158
$ws.Range("M12").Value = 18.92932567706474
$ws.Range("C13").Value = 4.900052222229425
$ws.Range("D13").Value = 8.194179456284933
$ws.Range("E13").Value = 12.71520496591781
$ws.Range("F13").Value = 36.17111211669978
$ws.Range("G13").Value = 40.77258208213338
$ws.Range("H13").Value = 17.03644081293636
$ws.Range("I13").Value = 25.33191177196179
$ws.Range("J13").Value = 9.814666818069483
$ws.Range("K13").Value = 18.14499037810255
$ws.Range("M13").Value = 18.9163404895793
$ws.Range("C14").Value = 4.888017686214505
$ws.Range("D14").Value = 8.189025822341266
$ws.Range("E14").Value = 12.71695216473094
$ws.Range("F14").Value = 36.17900265134369
$ws.Range("G14").Value = 40.7746118691752
$ws.Range("H14").Value = 17.0471433393674
$ws.Range("I14").Value = 25.33967255026884
$ws.Range("J14").Value = 9.819725900035516
$ws.Range("K14").Value = 18.05234747226523
$ws.Range("M14").Value = 18.87399134598015
$ws.Range("C15").Value = 4.880624752603897
$ws.Range("D15").Value = 8.185880081280597
$ws.Range("E15").Value = 12.7180595090108
$ws.Range("F15").Value = 36.18406638869724
$ws.Range("G15").Value = 40.7762021864609
$ws.Range("H15").Value = 17.0537799926281
$ws.Range("I15").Value = 25.34458556725394
$ws.Range("J15").Value = 9.822847936732023
$ws.Range("K15").Value = 17.99539196691952
$ws.Range("M15").Value = 18.84805085355563
$ws.Range("C16").Value = 4.837936574254361
$ws.Range("D16").Value = 8.168024126785884
$ws.Range("E16").Value = 12.72497463581158
$ws.Range("F16").Value = 36.21661108224362
$ws.Range("G16").Value = 40.79062168496581
$ws.Range("H16").Value = 17.09305184426742
$ws.Range("I16").Value = 25.37517626361165
$ws.Range("J16").Value = 9.841094937006661
$ws.Range("K16").Value = 17.6658794817418
$ws.Range("M16").Value = 18.69942434528057
$ws.Range("C17").Value = 4.811468090554337
$ws.Range("D17").Value = 8.157225138735255
$ws.Range("E17").Value = 12.72972512477125
$ws.Range("F17").Value = 36.23971460377248
$ws.Range("G17").Value = 40.80418749674731
$ws.Range("H17").Value = 17.11824765467353
$ws.Range("I17").Value = 25.3961101309336
$ws.Range("J17").Value = 9.852606685617268
$ws.Range("K17").Value = 17.461038195777
$ws.Range("M17").Value = 18.60830934730537
$ws.Range("C18").Value = 4.796143750886503
$ws.Range("D18").Value = 8.151071796996966
$ws.Range("E18").Value = 12.73264434742073
$ws.Range("F18").Value = 36.25415305786064
$ws.Range("G18").Value = 40.81371852932051
$ws.Range("H18").Value = 17.1331442442377
$ws.Range("I18").Value = 25.40894466698188
$ws.Range("J18").Value = 9.85934483356303
$ws.Range("K18").Value = 17.34226568972562
$ws.Range("M18").Value = 18.55594037760781
$ws.Range("C19").Value = 4.790938221952295
$ws.Range("D19").Value = 8.148998424923946
$ws.Range("E19").Value = 12.73366484695297
$ws.Range("F19").Value = 36.25923883212719
$ws.Range("G19").Value = 40.81724179330855
$ws.Range("H19").Value = 17.13825737181601
$ws.Range("I19").Value = 25.41342632433785
$ws.Range("J19").Value = 9.861646341269841
$ws.Range("K19").Value = 17.30189134579463
$ws.Range("M19").Value = 18.53821709214794
$ws.Range("C20").Value = 4.814296156849137
$ws.Range("D20").Value = 8.158368733408114
$ws.Range("E20").Value = 12.72920008862405
$ws.Range("F20").Value = 36.237136102412
$ws.Range("G20").Value = 40.8025643703714
$ws.Range("H20").Value = 17.11552361737679
$ws.Range("I20").Value = 25.39379945625506
$ws.Range("J20").Value = 9.851369144490119
$ws.Range("K20").Value = 17.48294342576381
$ws.Range("M20").Value = 18.61800510377495
$ws.Range("C21").Value = 4.891557002542068
$ws.Range("D21").Value = 8.190537280887208
$ws.Range("E21").Value = 12.71643123181443
$ws.Range("F21").Value = 36.17663694499618
$ws.Range("G21").Value = 40.77394332473865
$ws.Range("H21").Value = 17.0439828522124
$ws.Range("I21").Value = 25.33735980679562
$ws.Range("J21").Value = 9.818235090113896
$ws.Range("K21").Value = 18.07960269401843
$ws.Range("M21").Value = 18.88643048540398
$ws.Range("C22").Value = 4.941099550592054
$ws.Range("D22").Value = 8.212056618369044
$ws.Range("E22").Value = 12.70974831585988
$ws.Range("F22").Value = 36.1474074125391
$ws.Range("G22").Value = 40.77075320930293
$ws.Range("H22").Value = 17.00085155425578
$ws.Range("I22").Value = 25.30759695443799
$ws.Range("J22").Value = 9.797618851279356
$ws.Range("K22").Value = 18.46029599154522
$ws.Range("M22").Value = 19.0618834453184
$ws.Range("C23").Value = 4.914747597882214
$ws.Range("D23").Value = 8.200527037897698
$ws.Range("E23").Value = 12.71316311506222
$ws.Range("F23").Value = 36.16206120506072
$ws.Range("G23").Value = 40.77103074110764
$ws.Range("H23").Value = 17.02353883505422
$ws.Range("I23").Value = 25.3228276457209
$ws.Range("J23").Value = 9.808527179668971
$ws.Range("K23").Value = 18.25799487877223
$ws.Range("M23").Value = 18.96825412035597
$ws.Range("C24").Value = 4.813017921497332
$ws.Range("D24").Value = 8.157851542295285
$ws.Range("E24").Value = 12.72943687152588
$ws.Range("F24").Value = 36.23829824289036
$ws.Range("G24").Value = 40.80329279244241
$ws.Range("H24").Value = 17.11675387398614
$ws.Range("I24").Value = 25.39484162277811
$ws.Range("J24").Value = 9.851928263346512
$ws.Range("K24").Value = 17.47304319752085
$ws.Range("M24").Value = 18.61362160429634
$ws.Range("C25").Value = 4.699486349150802
$ws.Range("D25").Value = 8.114045932914703
$ws.Range("E25").Value = 12.75413535074767
$ws.Range("F25").Value = 36.36444315203056
$ws.Range("G25").Value = 40.90406637710417
$ws.Range("H25").Value = 17.23275821184457
$ws.Range("I25").Value = 25.50283568547118
$ws.Range("J25").Value = 9.903213935334399
$ws.Range("K25").Value = 16.59047600354477
$ws.Range("M25").Value = 18.232736007618
